# ---------------------------------------------------------------------
# Appends the 28th/29th.../33rd "day" blocks of class notes to the
# FSDS-10AM sheet (rows 326-378): a date + day-number on the first row
# of each day, followed by several note lines in column C only.
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSDS-10AM")

# Reuse the number formats already used lower in the sheet (row 312 is
# the previous day-header row) so the new cells land on the very same
# style indices instead of Excel minting brand-new ones.
$dateRef = $ws.Cells.Item(312, 1)
$numRef  = $ws.Cells.Item(312, 2)
$textRef = $ws.Cells.Item(312, 3)

$ws.Cells.Item(326, 1).Value = 45132
$dateRef.Copy()
$ws.Cells.Item(326, 1).PasteSpecial(-4122)
$ws.Cells.Item(326, 2).Value = 24
$numRef.Copy()
$ws.Cells.Item(326, 2).PasteSpecial(-4122)
$ws.Cells.Item(326, 3).Value = "new dataset project - movie rating anlaysis using pandas & seaborn"
$textRef.Copy()
$ws.Cells.Item(326, 3).PasteSpecial(-4122)

$ws.Cells.Item(327, 3).Value = ".astype to convert data type to other "
$textRef.Copy()
$ws.Cells.Item(327, 3).PasteSpecial(-4122)

$ws.Cells.Item(328, 3).Value = "discussed what is the problem statement , use case of the project"
$textRef.Copy()
$ws.Cells.Item(328, 3).PasteSpecial(-4122)

$ws.Cells.Item(329, 3).Value = "being a dataanalyst what we need to do . What is the goal & we are working on suggestion "
$textRef.Copy()
$ws.Cells.Item(329, 3).PasteSpecial(-4122)

$ws.Cells.Item(331, 3).Value = "distplot | histogram"
$textRef.Copy()
$ws.Cells.Item(331, 3).PasteSpecial(-4122)

$ws.Cells.Item(332, 3).Value = "uniform distribution "
$textRef.Copy()
$ws.Cells.Item(332, 3).PasteSpecial(-4122)

$ws.Cells.Item(333, 3).Value = "normal distribution | binomial || bell curve | 0 skewness | 0 syeemter | gausiaan "
$textRef.Copy()
$ws.Cells.Item(333, 3).PasteSpecial(-4122)

$ws.Cells.Item(334, 3).Value = "nomal distribution | mean = medain = mode"
$textRef.Copy()
$ws.Cells.Item(334, 3).PasteSpecial(-4122)

$ws.Cells.Item(335, 3).Value = "how to imputer hue in the graph "
$textRef.Copy()
$ws.Cells.Item(335, 3).PasteSpecial(-4122)

$ws.Cells.Item(330, 3).Value = "jointplot || kine - reg, scatter, kde, hex,resid"
$textRef.Copy()
$ws.Cells.Item(330, 3).PasteSpecial(-4122)

$ws.Cells.Item(336, 1).Value = 45133
$dateRef.Copy()
$ws.Cells.Item(336, 1).PasteSpecial(-4122)
$ws.Cells.Item(336, 2).Value = 25
$numRef.Copy()
$ws.Cells.Item(336, 2).PasteSpecial(-4122)
$ws.Cells.Item(336, 3).Value = "Final discussion what we learn so far - "
$textRef.Copy()
$ws.Cells.Item(336, 3).PasteSpecial(-4122)

$ws.Cells.Item(337, 3).Value = "1> category datatype in python"
$textRef.Copy()
$ws.Cells.Item(337, 3).PasteSpecial(-4122)

$ws.Cells.Item(338, 3).Value = "2> jointplots"
$textRef.Copy()
$ws.Cells.Item(338, 3).PasteSpecial(-4122)

$ws.Cells.Item(339, 3).Value = "3> histogram"
$textRef.Copy()
$ws.Cells.Item(339, 3).PasteSpecial(-4122)

$ws.Cells.Item(340, 3).Value = "4> stacked histograms"
$textRef.Copy()
$ws.Cells.Item(340, 3).PasteSpecial(-4122)

$ws.Cells.Item(341, 3).Value = "5> Kde plot"
$textRef.Copy()
$ws.Cells.Item(341, 3).PasteSpecial(-4122)

$ws.Cells.Item(342, 3).Value = "6> subplot"
$textRef.Copy()
$ws.Cells.Item(342, 3).PasteSpecial(-4122)

$ws.Cells.Item(343, 3).Value = "7> violin plots"
$textRef.Copy()
$ws.Cells.Item(343, 3).PasteSpecial(-4122)

$ws.Cells.Item(344, 3).Value = "8> Factet grid"
$textRef.Copy()
$ws.Cells.Item(344, 3).PasteSpecial(-4122)

$ws.Cells.Item(345, 3).Value = "9> Building dashboards"
$textRef.Copy()
$ws.Cells.Item(345, 3).PasteSpecial(-4122)

$ws.Cells.Item(346, 3).Value = "IRIS FLOWER - SEPAL LENGTH|SEPAL WIDTH | PETAL LENTH | PETEAL WIDTH"
$textRef.Copy()
$ws.Cells.Item(346, 3).PasteSpecial(-4122)

$ws.Cells.Item(347, 3).Value = "SEPAL > PETAL "
$textRef.Copy()
$ws.Cells.Item(347, 3).PasteSpecial(-4122)

$ws.Cells.Item(348, 3).Value = "PROJECT TO COMPLETE IRIS DATA ANALYSIS "
$textRef.Copy()
$ws.Cells.Item(348, 3).PasteSpecial(-4122)

$ws.Cells.Item(349, 1).Value = 45134
$dateRef.Copy()
$ws.Cells.Item(349, 1).PasteSpecial(-4122)
$ws.Cells.Item(349, 2).Value = 26
$numRef.Copy()
$ws.Cells.Item(349, 2).PasteSpecial(-4122)
$ws.Cells.Item(349, 3).Value = "discussed resume project "
$textRef.Copy()
$ws.Cells.Item(349, 3).PasteSpecial(-4122)

$ws.Cells.Item(350, 3).Value = "bank loan default risk anlaysis "
$textRef.Copy()
$ws.Cells.Item(350, 3).PasteSpecial(-4122)

$ws.Cells.Item(351, 3).Value = "discussed interview based question. "
$textRef.Copy()
$ws.Cells.Item(351, 3).PasteSpecial(-4122)

$ws.Cells.Item(352, 3).Value = "how to analysis , how to explore the dataset "
$textRef.Copy()
$ws.Cells.Item(352, 3).PasteSpecial(-4122)

$ws.Cells.Item(353, 3).Value = "I shared  kaggle link to practise, this is resume project "
$textRef.Copy()
$ws.Cells.Item(353, 3).PasteSpecial(-4122)

$ws.Cells.Item(354, 3).Value = "everyone to need to complete the project to have good exposure on data anlaytics skillls"
$textRef.Copy()
$ws.Cells.Item(354, 3).PasteSpecial(-4122)

$ws.Cells.Item(355, 1).Value = 45135
$dateRef.Copy()
$ws.Cells.Item(355, 1).PasteSpecial(-4122)
$ws.Cells.Item(355, 2).Value = 27
$numRef.Copy()
$ws.Cells.Item(355, 2).PasteSpecial(-4122)
$ws.Cells.Item(355, 3).Value = "how to find out which attribute are highly corelate with dv ==>"
$textRef.Copy()
$ws.Cells.Item(355, 3).PasteSpecial(-4122)

$ws.Cells.Item(356, 3).Value = "correlation = df.corr()"
$textRef.Copy()
$ws.Cells.Item(356, 3).PasteSpecial(-4122)

$ws.Cells.Item(357, 3).Value = "correlation['target'].sort_values(ascending=False)"
$textRef.Copy()
$ws.Cells.Item(357, 3).PasteSpecial(-4122)

$ws.Cells.Item(358, 3).Value = "eda project -- heart diseas analysis "
$textRef.Copy()
$ws.Cells.Item(358, 3).PasteSpecial(-4122)

$ws.Cells.Item(359, 3).Value = "seaborn  project- fifa dataset analysis "
$textRef.Copy()
$ws.Cells.Item(359, 3).PasteSpecial(-4122)

$ws.Cells.Item(360, 3).Value = "business anlaysis "
$textRef.Copy()
$ws.Cells.Item(360, 3).PasteSpecial(-4122)

$ws.Cells.Item(361, 3).Value = "even though business domains are difference methods are same algorithms are constant "
$textRef.Copy()
$ws.Cells.Item(361, 3).PasteSpecial(-4122)

$ws.Cells.Item(362, 3).Value = "business analys or data analyst or data scient "
$textRef.Copy()
$ws.Cells.Item(362, 3).PasteSpecial(-4122)

$ws.Cells.Item(363, 3).Value = "1st imp point understanding attribute or column name or feature "
$textRef.Copy()
$ws.Cells.Item(363, 3).PasteSpecial(-4122)

$ws.Cells.Item(364, 3).Value = "2nd most imp point is understand business well & dependent variabe"
$textRef.Copy()
$ws.Cells.Item(364, 3).PasteSpecial(-4122)

$ws.Cells.Item(365, 3).Value = "3rd point is based on dependent variable we select perfect algorithm"
$textRef.Copy()
$ws.Cells.Item(365, 3).PasteSpecial(-4122)

$ws.Cells.Item(366, 3).Value = "how to choose right machine learning algorithm"
$textRef.Copy()
$ws.Cells.Item(366, 3).PasteSpecial(-4122)

$ws.Cells.Item(367, 3).Value = "starting introduced about db & skewness we will talk later "
$textRef.Copy()
$ws.Cells.Item(367, 3).PasteSpecial(-4122)

$ws.Cells.Item(368, 1).Value = 45136
$dateRef.Copy()
$ws.Cells.Item(368, 1).PasteSpecial(-4122)
$ws.Cells.Item(368, 2).Value = 28
$numRef.Copy()
$ws.Cells.Item(368, 2).PasteSpecial(-4122)
$ws.Cells.Item(368, 3).Value = "sql introducetion "
$textRef.Copy()
$ws.Cells.Item(368, 3).PasteSpecial(-4122)

$ws.Cells.Item(369, 3).Value = "server - collection of db "
$textRef.Copy()
$ws.Cells.Item(369, 3).PasteSpecial(-4122)

$ws.Cells.Item(370, 3).Value = "db - collection of schma"
$textRef.Copy()
$ws.Cells.Item(370, 3).PasteSpecial(-4122)

$ws.Cells.Item(371, 3).Value = "schema - collection of tables"
$textRef.Copy()
$ws.Cells.Item(371, 3).PasteSpecial(-4122)

$ws.Cells.Item(372, 3).Value = "tables - collection of data types "
$textRef.Copy()
$ws.Cells.Item(372, 3).PasteSpecial(-4122)

$ws.Cells.Item(373, 3).Value = "datatype - int | char | varchar | logical "
$textRef.Copy()
$ws.Cells.Item(373, 3).PasteSpecial(-4122)

$ws.Cells.Item(374, 3).Value = "export raw data from the database using sql querying "
$textRef.Copy()
$ws.Cells.Item(374, 3).PasteSpecial(-4122)

$ws.Cells.Item(375, 3).Value = "steps to export data from db"
$textRef.Copy()
$ws.Cells.Item(375, 3).PasteSpecial(-4122)

$ws.Cells.Item(376, 3).Value = "team please please work on the project data extraction using sq & python "
$textRef.Copy()
$ws.Cells.Item(376, 3).PasteSpecial(-4122)

$ws.Cells.Item(377, 3).Value = "prepare ppt for data extraction steps "
$textRef.Copy()
$ws.Cells.Item(377, 3).PasteSpecial(-4122)

$ws.Cells.Item(378, 3).Value = "discused about website reference --> w3school python & w3school in sql "
$textRef.Copy()
$ws.Cells.Item(378, 3).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the author's final cursor position/selection after the paste.
[void]$ws.Activate()
[void]$ws.Range("C369").Select()

Write-Host "Added rows 326:378 to FSDS-10AM"